# Daily "tick" update for the shop water-delivery tracker.
#
# Columns: A=行号 B=店铺名称 C=地址 D=总天(total days) E=剩余(days remaining)
#          F=开始时间(start date, yyyymmdd)  G/H/I = 备注(notes)
#
# For every data row, the remaining-days count (E) is recomputed from the
# start date (F) and the total-day allotment (D) against "today". When the
# countdown would hit zero (or has already passed), the cycle restarts:
# the start date (F) is reset to today and the remaining days (E) reset to
# the full allotment (D). Rows whose start date can't be parsed are left
# untouched (matches source data quirks, e.g. a mistyped date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$invCulture = [System.Globalization.CultureInfo]::InvariantCulture

# "Today" for this automated run.
$today = Get-Date -Year 2025 -Month 12 -Day 29
$todaySerial = [math]::Floor($today.ToOADate())
$todayYmd = [int]$today.ToString("yyyyMMdd")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $dVal -or $null -eq $fVal) {
        continue
    }

    $totalDays = [int]$dVal
    $startYmd = [string][int64]$fVal

    try {
        $startDate = [datetime]::ParseExact($startYmd, "yyyyMMdd", $invCulture)
    } catch {
        # Unparseable start date (bad source data) - leave row untouched.
        continue
    }

    $startSerial = [math]::Floor($startDate.ToOADate())
    $endSerial = $startSerial + $totalDays
    $remaining = $endSerial - $todaySerial

    if ($remaining -le 0) {
        $ws.Cells.Item($r, 5).Value2 = $totalDays
        $ws.Cells.Item($r, 6).Value2 = $todayYmd
    } else {
        $ws.Cells.Item($r, 5).Value2 = $remaining
    }
}
